$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 200, shifting old rows 200-301 down to 202-303.
$ws.Range("A200:A201").EntireRow.Insert()

# New row 200
$ws.Cells.Item(200, 1).Value = 10
$ws.Cells.Item(200, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(200, 3).Value = "La Araucanía"
$ws.Cells.Item(200, 4).Value = "2021-11-23"
$ws.Cells.Item(200, 5).Value = 9
$ws.Cells.Item(200, 6).Value = "Fruta"
$ws.Cells.Item(200, 7).Value = 100103
$ws.Cells.Item(200, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(200, 9).Value = 100103006
$ws.Cells.Item(200, 10).Value = "Nectarín"
$ws.Cells.Item(200, 11).Value = "Early Glo"
$ws.Cells.Item(200, 12).Value = "Especial"
$ws.Cells.Item(200, 13).Value = 125
$ws.Cells.Item(200, 14).Value = 36000
$ws.Cells.Item(200, 15).Value = 36000
$ws.Cells.Item(200, 16).Value = 36000
$ws.Cells.Item(200, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item(200, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(200, 19).Value = 2000
$ws.Cells.Item(200, 20).Value = 18

# New row 201
$ws.Cells.Item(201, 1).Value = 10
$ws.Cells.Item(201, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(201, 3).Value = "La Araucanía"
$ws.Cells.Item(201, 4).Value = "2021-11-23"
$ws.Cells.Item(201, 5).Value = 9
$ws.Cells.Item(201, 6).Value = "Fruta"
$ws.Cells.Item(201, 7).Value = 100103
$ws.Cells.Item(201, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(201, 9).Value = 100103006
$ws.Cells.Item(201, 10).Value = "Nectarín"
$ws.Cells.Item(201, 11).Value = "Early Glo"
$ws.Cells.Item(201, 12).Value = "Primera"
$ws.Cells.Item(201, 13).Value = 260
$ws.Cells.Item(201, 14).Value = 28000
$ws.Cells.Item(201, 15).Value = 28000
$ws.Cells.Item(201, 16).Value = 28000
$ws.Cells.Item(201, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item(201, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(201, 19).Value = 1556
$ws.Cells.Item(201, 20).Value = 18
